$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-18 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-19 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("415×2=", $true, $false, $false, $false, $false, $true, 1, $false, "720×4=", 2) | Out-Null
$d.Content.Find.Execute("489×4=", $true, $false, $false, $false, $false, $true, 1, $false, "921×2=", 2) | Out-Null
$d.Content.Find.Execute("520×3=", $true, $false, $false, $false, $false, $true, 1, $false, "113×9=", 2) | Out-Null
$d.Content.Find.Execute("882×3=", $true, $false, $false, $false, $false, $true, 1, $false, "563×5=", 2) | Out-Null
$d.Content.Find.Execute("251×6=", $true, $false, $false, $false, $false, $true, 1, $false, "930×3=", 2) | Out-Null
$d.Content.Find.Execute("613×2=", $true, $false, $false, $false, $false, $true, 1, $false, "469×7=", 2) | Out-Null
$d.Content.Find.Execute("869×9=", $true, $false, $false, $false, $false, $true, 1, $false, "747×5=", 2) | Out-Null
$d.Content.Find.Execute("811×7=", $true, $false, $false, $false, $false, $true, 1, $false, "280×5=", 2) | Out-Null
$d.Content.Find.Execute("135×3=", $true, $false, $false, $false, $false, $true, 1, $false, "523×5=", 2) | Out-Null
$d.Content.Find.Execute("527×6=", $true, $false, $false, $false, $false, $true, 1, $false, "484×5=", 2) | Out-Null
$d.Content.Find.Execute("772×8=", $true, $false, $false, $false, $false, $true, 1, $false, "756×7=", 2) | Out-Null
$d.Content.Find.Execute("448×8=", $true, $false, $false, $false, $false, $true, 1, $false, "331×6=", 2) | Out-Null
$d.Content.Find.Execute("146×5=", $true, $false, $false, $false, $false, $true, 1, $false, "536×3=", 2) | Out-Null
$d.Content.Find.Execute("909×2=", $true, $false, $false, $false, $false, $true, 1, $false, "730×6=", 2) | Out-Null
$d.Content.Find.Execute("624×4=", $true, $false, $false, $false, $false, $true, 1, $false, "910×9=", 2) | Out-Null
$d.Content.Find.Execute("434×5=", $true, $false, $false, $false, $false, $true, 1, $false, "406×5=", 2) | Out-Null
$d.Content.Find.Execute("633×2=", $true, $false, $false, $false, $false, $true, 1, $false, "400×5=", 2) | Out-Null
$d.Content.Find.Execute("253×4=", $true, $false, $false, $false, $false, $true, 1, $false, "172×9=", 2) | Out-Null
$d.Content.Find.Execute("955×8=", $true, $false, $false, $false, $false, $true, 1, $false, "122×4=", 2) | Out-Null
$d.Content.Find.Execute("176×3=", $true, $false, $false, $false, $false, $true, 1, $false, "389×3=", 2) | Out-Null
$d.Content.Find.Execute("342×2=", $true, $false, $false, $false, $false, $true, 1, $false, "495×7=", 2) | Out-Null
$d.Content.Find.Execute("173×2=", $true, $false, $false, $false, $false, $true, 1, $false, "716×9=", 2) | Out-Null
$d.Content.Find.Execute("675×7=", $true, $false, $false, $false, $false, $true, 1, $false, "696×3=", 2) | Out-Null
$d.Content.Find.Execute("582×6=", $true, $false, $false, $false, $false, $true, 1, $false, "596×3=", 2) | Out-Null
$d.Content.Find.Execute("962×5=", $true, $false, $false, $false, $false, $true, 1, $false, "911×3=", 2) | Out-Null

Write-Output "Replacements applied"
